{"js": "// Office.js (Word JavaScript API) script.\n// Applies the edits described by the diff:\n//  1. Fix \"events modle\" -> \"events model\", \"delete a upcoming\" -> \"delete an upcoming\",\n//     and extend \"...for their resturant. \" -> \"...for their restaurant like bands and\n//     other shows happening. \", moving the \"_GoBack\" bookmark to right after the newly\n//     typed sentence (this is where Word leaves _GoBack after the last edit).\n//  2. Change the highlighted shell snippet \"module events\" -> \"module upcoming-events\".\n//  3. Normalize (merge into single runs) the \"endTime...\" and \"Details...\" paragraphs\n//     (their text does not otherwise change), removing the stale \"_GoBack\" bookmark that\n//     used to live in the \"Details...\" paragraph.\n//  4. Remove the stray empty paragraph between the \"Cover\" and \"Created\" bullet lines.\n//  5. Remove the duplicated \"View\"/\"Controllers\" block (13 paragraphs) that had been\n//     accidentally left in twice.\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replaceText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly 1 match for \" + JSON.stringify(searchText) + \" but found \" + results.items.length);\n  }\n  results.items[0].insertText(replaceText, \"Replace\");\n  await context.sync();\n}\n\n// --- Edit 1: events paragraph -------------------------------------------------\n// Drop the old \"_GoBack\" bookmark first (it currently sits in the \"Details...\" paragraph);\n// Word will re-create it at the new last-edited location below.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nawait replaceOnce(\"events modle\", \"events model\");\nawait replaceOnce(\"delete a upcoming\", \"delete an upcoming\");\nawait replaceOnce(\n  \"their resturant. \",\n  \"their restaurant like bands and other shows happening. \"\n);\n\n// Place a collapsed \"_GoBack\" bookmark right after \"...happening\" and before the\n// trailing \". \" -- matching where Word drops it after the final keystroke.\n{\n  const marker = body.search(\"restaurant like bands and other shows happening\", { matchCase: true });\n  marker.load(\"items\");\n  await context.sync();\n  const endRange = marker.items[0].getRange(\"End\");\n  endRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- Edit 2: \"module events\" -> \"module upcoming-events\" (keeps yellow highlight) ---\nawait replaceOnce(\"module events\", \"module upcoming-events\");\n\n// --- Edit 3: normalize runs (no text change) on the endTime/Details paragraphs ------\n{\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n\n  let endTimeParagraph = null;\n  let detailsParagraph = null;\n  for (const p of paragraphs.items) {\n    if (p.text === \"endTime \u2013 the time the event will end [String]\") {\n      endTimeParagraph = p;\n    } else if (p.text.indexOf(\"Details \u2013 details about the event\") === 0) {\n      detailsParagraph = p;\n    }\n  }\n\n  if (endTimeParagraph) {\n    const t = endTimeParagraph.text;\n    endTimeParagraph.getRange(\"Whole\").insertText(t, \"Replace\");\n    await context.sync();\n  }\n  if (detailsParagraph) {\n    const t = detailsParagraph.text;\n    detailsParagraph.getRange(\"Whole\").insertText(t, \"Replace\");\n    await context.sync();\n  }\n}\n\n// --- Edit 4 & 5: remove stray empty paragraph + duplicated View/Controllers block ---\n{\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n\n  const items = paragraphs.items;\n\n  // Find the \"Cover\" paragraph to anchor the single stray empty paragraph after it.\n  let coverIndex = -1;\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(\"Cover \u2013 the price of the menu item\") === 0) {\n      coverIndex = i;\n      break;\n    }\n  }\n  if (coverIndex === -1) {\n    throw new Error(\"Could not find the 'Cover' paragraph\");\n  }\n  const strayEmptyIndex = coverIndex + 1;\n  if (items[strayEmptyIndex].text !== \"\") {\n    throw new Error(\"Expected an empty paragraph after 'Cover' but found \" + JSON.stringify(items[strayEmptyIndex].text));\n  }\n\n  // Find the \"User \u2013 the user who added the menu item\" paragraph that follows \"Cover\"\n  // (that text also appears once earlier, for the unrelated \"Menuitem\" section, so we\n  // must anchor the search after \"Cover\" to land on the right occurrence). The\n  // duplicated View/Controllers block runs from just after it through the second copy\n  // of \"<module name>.client.controller.js\", followed by one more empty paragraph.\n  let userIndex = -1;\n  for (let i = coverIndex + 1; i < items.length; i++) {\n    if (items[i].text.indexOf(\"User \u2013 the user who added the menu item\") === 0) {\n      userIndex = i;\n      break;\n    }\n  }\n  if (userIndex === -1) {\n    throw new Error(\"Could not find the 'User' paragraph\");\n  }\n\n  // Walk forward from just after \"User...\" to find the second\n  // \"<module name>.client.controller.js\" occurrence; the block to delete ends at the\n  // following empty paragraph (inclusive).\n  let controllerJsIndex = -1;\n  for (let i = userIndex + 1; i < items.length; i++) {\n    if (items[i].text === \"<module name>.client.controller.js\") {\n      controllerJsIndex = i;\n      break;\n    }\n  }\n  if (controllerJsIndex === -1) {\n    throw new Error(\"Could not find the duplicated '<module name>.client.controller.js' paragraph\");\n  }\n  const blockStart = userIndex + 1; // first paragraph of the duplicated block (empty paragraph)\n  const blockEnd = controllerJsIndex + 1; // trailing empty paragraph right after the duplicated list item\n\n  if (items[blockEnd].text !== \"\") {\n    throw new Error(\"Expected an empty paragraph after the duplicated controller.js line but found \" + JSON.stringify(items[blockEnd].text));\n  }\n\n  // Delete from the bottom up so earlier indices stay valid.\n  for (let i = blockEnd; i >= blockStart; i--) {\n    items[i].delete();\n  }\n  items[strayEmptyIndex].delete();\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the edits described by the diff:\n#  1. Fix \"events modle\" -> \"events model\", \"delete a upcoming\" -> \"delete an upcoming\",\n#     and extend \"...for their resturant. \" -> \"...for their restaurant like bands and\n#     other shows happening. \", moving the \"_GoBack\" bookmark to right after the newly\n#     typed sentence (this is where Word leaves _GoBack after the last edit).\n#  2. Change the highlighted shell snippet \"module events\" -> \"module upcoming-events\".\n#  3. Normalize (merge into single runs) the \"endTime...\" and \"Details...\" paragraphs\n#     (their text does not otherwise change), removing the stale \"_GoBack\" bookmark that\n#     used to live in the \"Details...\" paragraph.\n#  4. Remove the stray empty paragraph between the \"Cover\" and \"Created\" bullet lines.\n#  5. Remove the duplicated \"View\"/\"Controllers\" block (13 paragraphs) that had been\n#     accidentally left in twice.\n\n$d = $word.ActiveDocument\n\n# wdReplace constants\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n# --- Edit 1: events paragraph -------------------------------------------------\n# Drop the old \"_GoBack\" bookmark first (it currently sits in the \"Details...\" paragraph);\n# we'll re-create it at the new last-edited location below.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$d.Content.Find.Execute(\"events modle\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"events model\", $wdReplaceOne) | Out-Null\n$d.Content.Find.Execute(\"delete a upcoming\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"delete an upcoming\", $wdReplaceOne) | Out-Null\n$d.Content.Find.Execute(\"resturant. \", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"restaurant like bands and other shows happening. \", $wdReplaceOne) | Out-Null\n\n# Place a collapsed \"_GoBack\" bookmark right after \"...happening\" and before the\n# trailing \". \" -- matching where Word drops it after the final keystroke.\n$marker = $d.Content\n$found = $marker.Find.Execute(\"restaurant like bands and other shows happening\")\nif ($found) {\n    $collapsePoint = $d.Range($marker.End, $marker.End)\n    $d.Bookmarks.Add(\"_GoBack\", $collapsePoint) | Out-Null\n}\n\n# --- Edit 2: \"module events\" -> \"module upcoming-events\" (keeps yellow highlight) ---\n$d.Content.Find.Execute(\"module events\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"module upcoming-events\", $wdReplaceOne) | Out-Null\n\n# --- Edit 3: normalize runs (no text change) on the endTime/Details paragraphs ------\n# NOTE: Paragraphs.Item(i).Range.Text includes the trailing paragraph-mark (\\r), so\n# exact-equality checks against literals must compare the trimmed text.\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.TrimEnd() -eq \"endTime \u2013 the time the event will end [String]\") {\n        $searchText = $t.TrimEnd()\n        $d.Content.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $searchText, $wdReplaceOne) | Out-Null\n        break\n    }\n}\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.StartsWith(\"Details \u2013 details about the event\")) {\n        $searchText = $t.TrimEnd()\n        $d.Content.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $searchText, $wdReplaceOne) | Out-Null\n        break\n    }\n}\n\n# --- Edit 4 & 5: remove stray empty paragraph + duplicated View/Controllers block ---\n$count = $d.Paragraphs.Count\n$coverIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.StartsWith(\"Cover \u2013 the price of the menu item\")) {\n        $coverIndex = $i\n        break\n    }\n}\nif ($coverIndex -eq -1) {\n    throw \"Could not find the 'Cover' paragraph\"\n}\n$strayEmptyIndex = $coverIndex + 1\nif ($d.Paragraphs.Item($strayEmptyIndex).Range.Text.TrimEnd() -ne \"\") {\n    throw \"Expected an empty paragraph after 'Cover'\"\n}\n\n# Find the \"User...\" paragraph that comes after \"Cover\" (that text also exists once\n# earlier, for the unrelated \"Menuitem\" section, so anchor the search after \"Cover\").\n$userIndex = -1\nfor ($i = $coverIndex + 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.StartsWith(\"User \u2013 the user who added the menu item\")) {\n        $userIndex = $i\n        break\n    }\n}\nif ($userIndex -eq -1) {\n    throw \"Could not find the 'User' paragraph\"\n}\n\n# Walk forward to find the second \"<module name>.client.controller.js\" occurrence;\n# the block to delete ends at the following empty paragraph (inclusive).\n$controllerJsIndex = -1\nfor ($i = $userIndex + 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq \"<module name>.client.controller.js\") {\n        $controllerJsIndex = $i\n        break\n    }\n}\nif ($controllerJsIndex -eq -1) {\n    throw \"Could not find the duplicated '<module name>.client.controller.js' paragraph\"\n}\n\n$blockStart = $userIndex + 1\n$blockEnd = $controllerJsIndex + 1\nif ($d.Paragraphs.Item($blockEnd).Range.Text.TrimEnd() -ne \"\") {\n    throw \"Expected an empty paragraph after the duplicated controller.js line\"\n}\n\n$firstPara = $d.Paragraphs.Item($blockStart)\n$lastPara = $d.Paragraphs.Item($blockEnd)\n$blockRange = $d.Range($firstPara.Range.Start, $lastPara.Range.End)\n$blockRange.Delete()\n\n$d.Paragraphs.Item($strayEmptyIndex).Range.Delete()\n"}
